# Enemy-Data.xlsx : "Enemies" sheet - rebalance HP values (nerf upgrades)
# Diff summary: column C (HP) lowered for several enemy types; column H
# ("Weg (Felder) bei 1 Schaden pro Sekunde") is a formula
# (=C*((D*40)/32)) so it recalculates automatically. Selection cell also
# moved from E18 to E17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enemies")

# Row 4 - Elite: HP 35 -> 28
$ws.Range("C4").Value = 28

# Row 5 - Tank: HP 100 -> 70
$ws.Range("C5").Value = 70

# Row 6 - Hover-Bike: HP 20 -> 18
$ws.Range("C6").Value = 18

# Row 7 - Jeep: HP 45 -> 30
$ws.Range("C7").Value = 30

# Row 8 - Space-Knight: HP 150 -> 90
$ws.Range("C8").Value = 90

# Row 9 - Roketeer: HP 20 -> 18
$ws.Range("C9").Value = 18

# Row 10 - Jet: HP 25 -> 18
$ws.Range("C10").Value = 18

# Row 11 - Ironclad: HP 250 -> 120
$ws.Range("C11").Value = 120

# Row 12 - Spec-Ops: HP 35 -> 30, and restyle C12 to match the plain
# bordered style used by the rest of the HP column (style index 8)
# instead of the row's fill-flagged style (index 12).
$ws.Range("C4").Copy()
$ws.Range("C12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C12").Value = 30

# Move the active-cell selection from E18 to E17 (matches the diff's
# <selection activeCell="E17" sqref="E17"/>)
$ws.Range("E17").Select()
